# Revert "Merge pull request #23 from AmreenNazia/Amreen"
# - sheet "Program": drop the two extra search-helper columns (D:E),
#   clear the stray "Devops"/"testing" scratch values, and rename the
#   "BDD" placeholder to "DA-course".
# - sheet "Class": left as-is (its values already match the reverted state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# B3 ("Devops") -> blank cell, same formatting as the rest of column B
$ws.Range("A5").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B3").Value = ""

# C3 ("BDD") -> "DA-course"
$ws.Range("C3").Value = "DA-course"

# A4 ("editprogram") -> blank cell, same formatting as the rest of column A
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4").Value = ""

# B4 ("Cybersecurity") -> blank cell, same formatting as the rest of column B
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B4").Value = ""

# C4 ("testing") -> remove entirely (row 4 only keeps A4/B4, like the other rows)
$ws.Range("C4").Clear()

# Drop the now-unused helper columns D (searchCreatedName/Devops) and
# E (SearchUpdatedName/Cybersecurity) entirely
$ws.Columns("D:E").Delete()

# Restore the original selection
[void]$ws.Range("A4").Select()
